# Update the Fitness column (C) values in the log worksheet.
# Rows 2-241 (generations 0-239) change from 7534/7295 to 7310.
# Rows 242-252 (generations 240-250) change from 7295 to 7293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C241").Value = 7310
$ws.Range("C242:C252").Value = 7293
